$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Get ready to play Cash Spin, an online slot game with modern features and retro-style design. Read our review today to experience the thrilling bonus rounds. Play now for free!</w:t></w:r></w:p>'
$metaRange.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of the doc
#    and rewrite the italic "meta description" paragraph that follows
#    it with the new AI image-generation prompt text.
# ---------------------------------------------------------------------
$dupTitle = "Play Cash Spin Free " + [char]0x2013 + " Read Our Slot Game Review"
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq $dupTitle) {
        $p.Range.Delete()
        break
    }
}

$oldBlurb = "Get ready to play Cash Spin, an online slot game with modern features and retro-style design. Read our review today to experience the thrilling bonus rounds. Play now for free!"
$newBlurb = 'Create a feature image for the game "Cash Spin" that fits the following criteria: - In cartoon style - Features a happy Maya warrior with glasses The image should be lively and colorful, with a cartoon-style depiction of a Maya warrior wearing a big smile and black-rimmed glasses. The warrior should be holding a wheel of fortune in one hand, indicating the game''s bonus feature, while holding a money bag in the other hand to represent the Coin Purse function. The background of the image can be a colorful mixture of the game''s classic symbols, such as diamonds, rubies, emeralds, and dollar signs, arranged in a fun and playful way. Overall, the image should be eye-catching and playful, reflecting the game''s fun and nostalgic atmosphere.'

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq $oldBlurb) {
        $r = $p.Range
        $inner = $d.Range($r.Start, $r.End - 1)
        $inner.Text = $newBlurb
        break
    }
}
